$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Adjustments to infiltration parameters (rows 14-17: Je, Jc, Jp, Jn)
$ws.Range("B14").Value = 0.01
$ws.Range("D14").Value = 0.01

$ws.Range("B15").Value = 0.01
$ws.Range("D15").Value = 0.01

$ws.Range("B16").Value = 0.01
$ws.Range("D16").Value = 0.01

$ws.Range("B17").Value = 0.01
$ws.Range("D17").Value = 0.01

# Update the active selection to match the saved state of the sheet
$ws.Range("G18").Select()
